# Test case updated 20-10
# RUNMANAGER: RetailLoginLogout / myAccountsSummary / myAccountsStatement test rows
# RETAIL_DATA: matching data rows + a new hyperlinked password cell

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("RETAIL_DATA")

# ---------------------------------------------------------------
# Sheet "RUNMANAGER"
# ---------------------------------------------------------------

# Row 2 : login/logout test -> retail login/logout test
$ws1.Range("A2").Value = "RetailLoginLogout"
$ws1.Range("B2").Value = "To Login into the application"

# Row 3 : newTest -> myAccountsSummary (and bump priority to 2)
$ws1.Range("A3").Value = "myAccountsSummary"
$ws1.Range("B3").Value = "To check Mouse Hover functionality on My Account menu"
$ws1.Range("D3").Value = "'2"

# Row 4 (new) : myAccountsStatement
$ws1.Range("A4").Value = "myAccountsStatement"
$ws1.Range("B4").Value = "To check account statement"
$ws1.Range("C4").Value = "no"
$ws1.Range("D4").Value = "'3"
$ws1.Range("E4").Value = "'1"

# Stray formatted-but-empty cells that ride along with the edited range
$ws1.Range("N3").Value = "'1"
$ws1.Range("N3").ClearContents() | Out-Null
$ws1.Range("O3").Value = "'1"
$ws1.Range("O3").ClearContents() | Out-Null
$ws1.Range("D5").Value = "'1"
$ws1.Range("D5").ClearContents() | Out-Null

# Column A got wider to fit the longer test-case names
$ws1.Columns.Item(1).ColumnWidth = 30.666666666666668

# ---------------------------------------------------------------
# Sheet "RETAIL_DATA"
# ---------------------------------------------------------------

# Row 3 : login/logout test -> retail login/logout test
$ws2.Range("A3").Value = "RetailLoginLogout"
$ws2.Range("B3").Value = "no"

# Row 4 : newTest -> myAccountsSummary
$ws2.Range("A4").Value = "myAccountsSummary"

# Row 5 (new) : myAccountsStatement, with a hyperlinked password cell
$ws2.Range("A5").Value = "myAccountsStatement"
$ws2.Range("B5").Value = "no"
$ws2.Range("C5").Value = "chrome"
$ws2.Range("D5").Value = "spcb"
$ws2.Range("E5").Value = "Asdf@123"
$ws2.Hyperlinks.Add($ws2.Range("E5"), "mailto:Asdf@123") | Out-Null
$ws2.Range("E5").Style = "Hyperlink"

# Column A widened (and now auto-fit) to match the longer names
$ws2.Columns.Item(1).ColumnWidth = 23.333333333333332

# ---------------------------------------------------------------
# Active sheet / selection bookkeeping
# (RETAIL_DATA was the active tab before; RUNMANAGER is now)
# ---------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("D9").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B21").Select() | Out-Null
